$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fix: "(in percents)" -> "(in percent)" in C2 ---
$ws.Range("C2").Value = "(in percent)"

# --- Row heights for rows 2 and 3 ---
$ws.Rows("2:3").RowHeight = 14.25

# --- Column widths ---
# A:C from 27.625 -> as close as COM rounding allows to 28.25
$ws.Columns("A:C").ColumnWidth = 27.57
# D:T from default(9) -> as close as COM rounding allows to 8.625
$ws.Columns("D:T").ColumnWidth = 7.8

# --- New column T (year 2023) ---
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 2.5
$ws.Range("T6").Value = 2.7
$ws.Range("T7").Value = 2.2000000000000002

# Copy number formats/styles from column S onto the new column T cells
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("S7").Copy()
$ws.Range("T7").PasteSpecial(-4122)
